$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Update the "Date" property (row 8, column B) to the new export timestamp. ---
$ws.Range("B8").Value = "2022-01-21T07:49:24+01:00"

# --- Make room for two extra "Contact" rows (one per IG author) right after the
# existing Contact row (row 11). Everything from the old row 12 ("Jurisdiction")
# down to the old row 23 ("Count") needs to move down by two rows.
#
# We deliberately avoid Rows.Insert() here: while it does shift rows down, this
# runtime's implementation always fabricates one brand-new (unused) cell style for
# the freshly inserted blank row, which would needlessly grow styles.xml. Instead
# we shift the data ourselves, row by row, from the bottom up (so a source row is
# always read before it gets overwritten).
for ($r = 23; $r -ge 12; $r--) {
    $destRow = $r + 2
    $srcRange = $ws.Range("A$r`:B$r")
    $dstRange = $ws.Range("A$destRow`:B$destRow")

    if ($destRow -gt 23) {
        # Rows 24/25 do not exist yet (original sheet only went to row 23). Paste
        # formats first so the brand-new cells actually pick up the source style;
        # otherwise a freshly-materialized cell ends up with no style at all.
        $srcRange.Copy()
        $dstRange.PasteSpecial(-4122)
    } else {
        # Destination row already has old content. Clear it first so that a blank
        # source cell actually leaves the destination blank instead of retaining
        # whatever stale value used to be there.
        $dstRange.ClearContents()
    }

    $srcRange.Copy()
    $dstRange.PasteSpecial(-4104)
}
$excel.CutCopyMode = 0

# --- Fill the two newly freed rows (12 and 13) with the new Contact entries. ---
$ws.Range("A12").Value = "Contact"
$ws.Range("B12").Value = "No display for ContactDetail"
$ws.Range("A13").Value = "Contact"
$ws.Range("B13").Value = "No display for ContactDetail"
